# Applies the cryptos list update described by the commit diff.
# Each changed cell is written as text (NumberFormat "@") so that
# numeric-looking strings (e.g. "0.999", "7.36") are preserved as
# literal text rather than being coerced to floating point numbers,
# matching the original inlineStr cell type. The style is reset back
# to "Normal" afterward so no stray style/number-format is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '62.509.42'
Set-TextValue 'E2' '  +2.68%  '
Set-TextValue 'D3' '2.969.27'
Set-TextValue 'E3' '  +2.07%  '
Set-TextValue 'D4' '0.999'
Set-TextValue 'E4' '  -0.15%  '
Set-TextValue 'D5' '597.25'
Set-TextValue 'D6' '144.82'
Set-TextValue 'E6' '  +0.11%  '
Set-TextValue 'D7' '0.998'
Set-TextValue 'D8' '2.970.07'
Set-TextValue 'E8' '  +2.12%  '
Set-TextValue 'D9' '0.504'
Set-TextValue 'E9' '  -0.12%  '
Set-TextValue 'D10' '7.36'
Set-TextValue 'E10' '  +6.73%  '
Set-TextValue 'D11' '0.144'
Set-TextValue 'E11' '  +2.50%  '
Set-TextValue 'D12' '0.447'
Set-TextValue 'E12' '  +1.90%  '
Set-TextValue 'D13' '0.0000232'
Set-TextValue 'E13' '  +3.34%  '
Set-TextValue 'D14' '33.70'
Set-TextValue 'E14' '  +0.99%  '
Set-TextValue 'E15' '  -0.06%  '
Set-TextValue 'D16' '3.453.94'
Set-TextValue 'E16' '  +1.84%  '
Set-TextValue 'D17' '62.310.57'
Set-TextValue 'E17' '  +2.44%  '
Set-TextValue 'B18' 'Polkadot'
Set-TextValue 'C18' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D18' '6.74'
Set-TextValue 'E18' '  +1.06%  '
Set-TextValue 'B19' 'WrappedEther'
Set-TextValue 'C19' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D19' '2.958.18'
Set-TextValue 'E19' '  +1.73%  '
Set-TextValue 'D20' '444.17'
Set-TextValue 'E20' '  +2.77%  '
Set-TextValue 'D21' '13.60'
Set-TextValue 'E21' '  +1.98%  '
Set-TextValue 'D22' '0.677'
Set-TextValue 'E22' '  +0.17%  '
Set-TextValue 'D23' '7.13'
Set-TextValue 'E23' '  +0.59%  '
Set-TextValue 'D24' '81.78'
Set-TextValue 'E24' '  +0.45%  '
Set-TextValue 'D25' '11.05'
Set-TextValue 'E25' '  +2.37%  '
Set-TextValue 'D26' '2.18'
Set-TextValue 'E26' '  -0.06%  '
Set-TextValue 'D27' '12.06'
Set-TextValue 'E27' '  +2.34%  '
Set-TextValue 'E28' '  -0.03%  '
Set-TextValue 'D29' '2.61'
Set-TextValue 'E29' '  +0.91%  '
Set-TextValue 'B30' 'NEARProtocol'
Set-TextValue 'C30' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D30' '7.02'
Set-TextValue 'E30' '  +1.10%  '
Set-TextValue 'B31' 'ImmutableX'
Set-TextValue 'C31' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D31' '2.13'
Set-TextValue 'E31' '  -6.21%  '
Set-TextValue 'D32' '26.52'
Set-TextValue 'E32' '  +0.15%  '
Set-TextValue 'B33' 'Hedera'
Set-TextValue 'C33' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D33' '0.108'
Set-TextValue 'E33' '  -0.90%  '
Set-TextValue 'B34' 'FirstDigitalUSD'
Set-TextValue 'C34' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D34' '1.00'
Set-TextValue 'E34' '  -0.14%  '
Set-TextValue 'D35' '0.0₃0885'
Set-TextValue 'E35' '  +3.16%  '
Set-TextValue 'D36' '1.00'
Set-TextValue 'E36' '  -0.50%  '
Set-TextValue 'D37' '5.64'
Set-TextValue 'E37' '  +0.60%  '
Set-TextValue 'D38' '49.64'
Set-TextValue 'E38' '  +0.14%  '
Set-TextValue 'B39' 'dogwifhat'
Set-TextValue 'C39' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D39' '2.94'
Set-TextValue 'E39' '  -1.71%  '
Set-TextValue 'B40' 'Stacks'
Set-TextValue 'C40' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D40' '2.01'
Set-TextValue 'E40' '  +1.89%  '
Set-TextValue 'D41' '8.63'
Set-TextValue 'E41' '  +1.16%  '
Set-TextValue 'D42' '0.120'
Set-TextValue 'E42' '  -0.89%  '
Set-TextValue 'D43' '0.282'
Set-TextValue 'E43' '  +0.48%  '
Set-TextValue 'D44' '39.87'
Set-TextValue 'E44' '  -2.75%  '
Set-TextValue 'D45' '2.693.52'
Set-TextValue 'E45' '  -0.13%  '
Set-TextValue 'B46' 'VeChain'
Set-TextValue 'C46' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D46' '0.0342'
Set-TextValue 'E46' '  -0.77%  '
Set-TextValue 'B47' 'Bittensor'
Set-TextValue 'C47' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D47' '367.28'
Set-TextValue 'E47' '  -2.28%  '
Set-TextValue 'B48' 'Monero'
Set-TextValue 'C48' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D48' '133.53'
Set-TextValue 'E48' '  +0.12%  '
Set-TextValue 'E49' '  +0.05%  '
Set-TextValue 'D50' '23.38'
Set-TextValue 'E50' '  -1.34%  '
Set-TextValue 'D51' '0.105'
Set-TextValue 'E51' '  -0.98%  '
